$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 492
$ws.Range("F5").Value = 2105
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 7931
$ws.Range("F8").Value = 267
$ws.Range("F9").Value = 52
$ws.Range("F11").Value = 232
$ws.Range("F12").Value = 1766
$ws.Range("F13").Value = 1540
$ws.Range("F14").Value = 1313
$ws.Range("F15").Value = 181
$ws.Range("F16").Value = 4000
$ws.Range("F17").Value = 695
$ws.Range("F18").Value = 26
$ws.Range("F19").Value = 1099
$ws.Range("F20").Value = 1230
$ws.Range("F21").Value = 427
$ws.Range("F22").Value = 6204
$ws.Range("F25").Value = 4211
$ws.Range("F26").Value = 704
$ws.Range("F27").Value = 1944
$ws.Range("F28").Value = 1166
$ws.Range("F29").Value = 301
$ws.Range("F30").Value = 1029
$ws.Range("F31").Value = 17
$ws.Range("F32").Value = 33
$ws.Range("F33").Value = 202
$ws.Range("F34").Value = 46
$ws.Range("F35").Value = 319
$ws.Range("F36").Value = 1151
$ws.Range("F37").Value = 502
$ws.Range("F38").Value = 1871
$ws.Range("F39").Value = 107
$ws.Range("F40").Value = 407
$ws.Range("F41").Value = 151
$ws.Range("F42").Value = 1138
$ws.Range("F43").Value = 555
$ws.Range("F45").Value = 32
$ws.Range("F48").Value = 170

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 990
$ws.Range("F5").Value = 120
$ws.Range("F11").Value = 678
$ws.Range("F12").Value = 369
$ws.Range("F13").Value = 404
$ws.Range("F15").Value = 204
$ws.Range("F16").Value = 107
$ws.Range("F20").Value = 172
$ws.Range("F22").Value = 83
$ws.Range("F24").Value = 228
$ws.Range("F25").Value = 93
$ws.Range("F30").Value = 270

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 453
$ws.Range("F6").Value = 1558
$ws.Range("F8").Value = 3098
$ws.Range("F9").Value = 929
$ws.Range("F10").Value = 1070
$ws.Range("F11").Value = 1261
$ws.Range("F12").Value = 1573

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 453
$ws.Range("F3").Value = 1558
$ws.Range("F5").Value = 492
$ws.Range("F7").Value = 3098
$ws.Range("F8").Value = 2105
$ws.Range("F9").Value = 7931
$ws.Range("F10").Value = 52
$ws.Range("F11").Value = 929
$ws.Range("F14").Value = 232
$ws.Range("F15").Value = 1766
$ws.Range("F16").Value = 1540
$ws.Range("F17").Value = 1261
$ws.Range("F18").Value = 1313
$ws.Range("F19").Value = 678
$ws.Range("F20").Value = 181
$ws.Range("F21").Value = 1573
$ws.Range("F22").Value = 4000
$ws.Range("F23").Value = 369
$ws.Range("F24").Value = 404
$ws.Range("F25").Value = 695
$ws.Range("F26").Value = 26
$ws.Range("F27").Value = 1099
$ws.Range("F28").Value = 1230
$ws.Range("F29").Value = 427
$ws.Range("F30").Value = 6204
$ws.Range("F32").Value = 704
$ws.Range("F33").Value = 1944
$ws.Range("F34").Value = 1166
$ws.Range("F35").Value = 301
$ws.Range("F36").Value = 33
$ws.Range("F37").Value = 172
$ws.Range("F38").Value = 202
$ws.Range("F39").Value = 83
$ws.Range("F40").Value = 502
$ws.Range("F41").Value = 1871
$ws.Range("F42").Value = 107
$ws.Range("F43").Value = 407
$ws.Range("F44").Value = 1138
$ws.Range("F46").Value = 555
$ws.Range("F47").Value = 270
$ws.Range("F49").Value = 170
